$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.130.95"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "'1.857.51"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'232.78"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.4692"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "'42.79"
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("D9").Value = "'0.2831"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D10").Value = "'0.06446"
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("D11").Value = "'20.71"
$ws.Range("E11").Value = "  -4.23%  "
$ws.Range("D12").Value = "'0.07687"
$ws.Range("E12").Value = "  -4.16%  "
$ws.Range("D13").Value = "'1.862.26"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").Value = "'93.34"
$ws.Range("E14").Value = "  -3.65%  "
$ws.Range("D15").Value = "'5.052"
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").Value = "'0.6780"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "'263.91"
$ws.Range("E17").Value = "  -2.00%  "
$ws.Range("D18").Value = "'30.106.48"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "'13.35"
$ws.Range("E19").Value = "  -4.41%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.000007511"
$ws.Range("E21").Value = "  -1.52%  "
$ws.Range("D22").Value = "'2.121.71"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "'5.144"
$ws.Range("E24").Value = "  -2.52%  "
$ws.Range("D25").Value = "'6.099"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("D26").Value = "'9.255"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("D27").Value = "'165.19"
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("D28").Value = "'18.43"
$ws.Range("E28").Value = "  -2.34%  "
$ws.Range("D29").Value = "'1.877"
$ws.Range("E29").Value = "  -3.61%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").Value = "'0.09805"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("D32").Value = "'1.446"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").Value = "'4.200"
$ws.Range("E33").Value = "  -3.90%  "
$ws.Range("D34").Value = "'3.975"
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("D35").Value = "'0.04646"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").Value = "'1.110"
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("D37").Value = "'0.6834"
$ws.Range("E37").Value = "  -2.33%  "
$ws.Range("D38").Value = "'2.714"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "'0.01805"
$ws.Range("E39").Value = "  -3.72%  "
$ws.Range("D40").Value = "'2.719"
$ws.Range("E40").Value = "  +3.52%  "
$ws.Range("D41").Value = "'6.281"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "'70.16"
$ws.Range("E42").Value = "  -2.76%  "
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "'0.8301"
$ws.Range("E44").Value = "  -1.30%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'102.33"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'1.873"
$ws.Range("E46").Value = "  -4.44%  "
$ws.Range("D47").Value = "'0.4023"
$ws.Range("E47").Value = "  -3.39%  "
$ws.Range("D48").Value = "'9.122"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").Value = "'921.23"
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("D50").Value = "'6.897"
$ws.Range("E50").Value = "  -2.37%  "
$ws.Range("D51").Value = "'34.04"
$ws.Range("E51").Value = "  -1.11%  "
